$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated CAD model / linked-parameter pose values (row 2)
$ws.Range("E2").Value = 0.500274
$ws.Range("F2").Value = -9.0432
$ws.Range("G2").Value = 3.69867
$ws.Range("I2").Value = -3.69083
$ws.Range("J2").Value = 3.6051
$ws.Range("K2").Value = 0.722841
$ws.Range("M2").Value = 11.0918

# row 3
$ws.Range("E3").Value = -4.99892
$ws.Range("F3").Value = 0.784414
$ws.Range("G3").Value = -0.922893
$ws.Range("I3").Value = 4.88145
$ws.Range("J3").Value = -1.55169
$ws.Range("K3").Value = -5.59543
$ws.Range("M3").Value = 11.9904

# row 4
$ws.Range("E4").Value = -9.96637
$ws.Range("F4").Value = -12.9539
$ws.Range("G4").Value = 14.0451
$ws.Range("I4").Value = -7.39793
$ws.Range("J4").Value = 0.538118
$ws.Range("K4").Value = 7.16814
$ws.Range("M4").Value = 15.4573

# row 5
$ws.Range("E5").Value = 15.7575
$ws.Range("F5").Value = -0.714453
$ws.Range("G5").Value = 4.79598
$ws.Range("I5").Value = 0.689914
$ws.Range("J5").Value = -4.21556
$ws.Range("K5").Value = 1.69537
$ws.Range("M5").Value = 18.8602

# row 6
$ws.Range("E6").Value = 3.12885
$ws.Range("F6").Value = -7.02329
$ws.Range("G6").Value = 14.568
$ws.Range("I6").Value = -2.85782
$ws.Range("J6").Value = -16.1458
$ws.Range("K6").Value = 14.1163
$ws.Range("M6").Value = 20.8349

# Update the active selection to match the updated pose block
$ws.Range("E2:M7").Select() | Out-Null
